# Updates prices/volumes/coin listings in the cryptos tracker sheet
# (mirrors the scheduled GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "94.737.11"
$ws.Range("E2").Value = "  +2.77%  "
$ws.Range("D3").Value = "3.108.23"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'237.23"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.28%  "
$ws.Range("D6").Value = "'612.86"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("E7").Value = "  +2.31%  "
$ws.Range("D8").Value = "'0.393"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").Value = "'1.00"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").Value = "'0.834"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +14.54%  "
$ws.Range("D11").Value = "3.102.85"
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("E12").Value = "  -2.02%  "
$ws.Range("E13").Value = "  -2.31%  "
$ws.Range("D14").Value = "94.296.77"
$ws.Range("E14").Value = "  +2.70%  "
$ws.Range("D15").Value = "'34.59"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.85%  "
$ws.Range("E16").Value = "  -2.60%  "
$ws.Range("D17").Value = "3.681.86"
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("D18").Value = "3.103.54"
$ws.Range("E18").Value = "  +0.50%  "
$ws.Range("E19").Value = "  +0.51%  "
$ws.Range("D20").Value = "'14.90"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.25%  "
$ws.Range("D21").Value = "'5.91"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.79%  "
$ws.Range("D22").Value = "'450.24"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.67%  "
$ws.Range("E23").Value = "  -1.49%  "
$ws.Range("D24").Value = "'8.96"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -3.80%  "
$ws.Range("D25").Value = "'8.31"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +5.52%  "
$ws.Range("D26").Value = "'5.62"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("B27").Value = "Litecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D27").Value = "'85.96"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.24%  "
$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D28").Value = "'12.12"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +3.62%  "
$ws.Range("D29").Value = "3.284.05"
$ws.Range("E29").Value = "  +0.56%  "
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("D31").Value = "'0.256"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +8.25%  "
$ws.Range("D32").Value = "'0.183"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +9.14%  "
$ws.Range("D33").Value = "'0.124"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -9.57%  "
$ws.Range("D34").Value = "'9.31"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.37%  "
$ws.Range("E35").Value = "  +0.27%  "
$ws.Range("D36").Value = "'7.89"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.30%  "
$ws.Range("E37").Value = "  +0.70%  "
$ws.Range("E38").Value = "  -0.69%  "
$ws.Range("D39").Value = "'1.90"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.97%  "
$ws.Range("D40").Value = "'0.455"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +5.20%  "
$ws.Range("E41").Value = "  +5.44%  "
$ws.Range("E42").Value = "  -1.13%  "
$ws.Range("D43").Value = "'471.43"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.77%  "
$ws.Range("D44").Value = "'3.68"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -13.48%  "
$ws.Range("D45").Value = "'3.23"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -5.12%  "
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").Value = "'160.07"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("D48").Value = "'0.690"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.34%  "
$ws.Range("D49").Value = "'1.86"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.15%  "
$ws.Range("D50").Value = "'4.41"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.28%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "'0.0324"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.73%  "
